$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold text-formatted numbers (e.g. "0.580", "1.60")
# that must stay text (not be coerced to numeric 0.58 / 1.6), matching the
# original workbook formatting, so force text format before assigning.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.906.97"
$ws.Range("E2").Value = "  -0.29%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.517.23"
$ws.Range("E3").Value = "  -1.23%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.98"
$ws.Range("E5").Value = "  -0.74%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.26"
$ws.Range("E6").Value = "  -2.61%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.515.50"
$ws.Range("E7").Value = "  -1.25%  "

# Row 9
$ws.Range("E9").Value = "  -1.03%  "

# Row 10
$ws.Range("E10").Value = "  +0.64%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.16"
$ws.Range("E11").Value = "  +2.65%  "

# Row 12
$ws.Range("E12").Value = "  -0.89%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.112.72"
$ws.Range("E13").Value = "  -1.34%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.64"
$ws.Range("E14").Value = "  +1.11%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000182"
$ws.Range("E15").Value = "  -1.00%  "

# Row 16
$ws.Range("E16").Value = "  +0.15%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.514.36"
$ws.Range("E17").Value = "  -1.54%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.829.93"
$ws.Range("E18").Value = "  -0.13%  "

# Row 19
$ws.Range("E19").Value = "  -1.29%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.34"
$ws.Range("E20").Value = "  -0.42%  "

# Row 21
$ws.Range("E21").Value = "  -2.79%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "392.34"
$ws.Range("E22").Value = "  +0.17%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.580"
$ws.Range("E23").Value = "  +0.11%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.658.51"
$ws.Range("E24").Value = "  -1.31%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.22"
$ws.Range("E25").Value = "  +0.04%  "

# Row 26
$ws.Range("E26").Value = "  +0.00%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000112"
$ws.Range("E27").Value = "  -4.51%  "

# Row 28
$ws.Range("B28").Value = "Fetch.AI"
$ws.Range("C28").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.60"
$ws.Range("E28").Value = "  +9.77%  "

# Row 29
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.70"
$ws.Range("E29").Value = "  +0.05%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -0.15%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.34"
$ws.Range("E32").Value = "  +0.23%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.521.28"
$ws.Range("E33").Value = "  -1.39%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.21"
$ws.Range("E34").Value = "  +0.65%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.145"
$ws.Range("E36").Value = "  +0.18%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.29"
$ws.Range("E37").Value = "  +5.53%  "

# Row 38
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.98"
$ws.Range("E38").Value = "  +1.01%  "

# Row 39
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.57"
$ws.Range("E39").Value = "  +1.09%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "168.57"
$ws.Range("E40").Value = "  -0.95%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0819"
$ws.Range("E41").Value = "  +1.03%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.823"
$ws.Range("E42").Value = "  -0.49%  "

# Row 43
$ws.Range("E43").Value = "  +2.32%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.70"
$ws.Range("E44").Value = "  -4.58%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.76"
$ws.Range("E45").Value = "  +0.32%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  -0.09%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.43"
$ws.Range("E47").Value = "  -0.95%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.65"
$ws.Range("E48").Value = "  -1.40%  "

# Row 49
$ws.Range("E49").Value = "  +0.13%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.388.81"
$ws.Range("E50").Value = "  -4.39%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.900"
$ws.Range("E51").Value = "  +4.38%  "
